# Update the "dSF" column (F) values for the specified rows to reflect
# repulled data / push all data / mean calculation corrections.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F4").Value = -1
$ws.Range("F7").Value = -7
$ws.Range("F9").Value = -5
$ws.Range("F13").Value = 0
